$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 becomes a blank cell (value removed)
$ws.Range("D3").ClearContents()

# Row 7: "Other" -> "Biogas", value updated
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 143.8823441355079

# New Row 8: "Other" (moved from row 7) with corrected value.
# Copy formatting from A7 (bold/border/centered style) onto A8 first.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 1442.378226233408
